$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header / account holder info
$ws.Range("C2").Value = "Hartmut"

# B3 holds a long digit-string that must remain TEXT (not get auto-converted
# to a number). Force text format, assign value, then restore the original
# look (General number format / same style) via a format-only copy from a
# neighboring cell that already uses the same style (s=8, General).
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "2570314725427075"
$ws.Range("C3").Copy()
$ws.Range("B3").PasteSpecial(-4122)

$ws.Range("C3").Value = "Mohaupt"

# Opening balance line
$ws.Range("D5").Value = "KONTOSTAND AM 29.04.2025"

# Transaction rows shift up (old row7-10 content moves into row6-9),
# and row 10 becomes an empty spacer row like row 11.

# Row 6
$ws.Range("B6").Value = "02.05."
$ws.Range("C6").Value = "03.05."
$ws.Range("D6").Value = "MITGLIEDSBEITRAG ZEUS BODYPOWER"
$ws.Range("E6").Value = "24,84-"

# Row 7
$ws.Range("B7").Value = "04.05."
$ws.Range("C7").Value = "05.05."
$ws.Range("D7").Value = "PAYPAL RBHGEW"
$ws.Range("E7").Value = "95,45-"

# Row 8
$ws.Range("B8").Value = "05.05."
$ws.Range("C8").Value = "06.05."
$ws.Range("D8").Value = "AMAZON.DE MKTPLC EU OQFJYI"
$ws.Range("E8").Value = "193,83-"

# Row 9
$ws.Range("B9").Value = "06.05."
$ws.Range("C9").Value = "07.05."
$ws.Range("D9").Value = "RECHNUNG VODAFONE GMBH 24840316"
$ws.Range("E9").Value = "39,05-"

# Row 10 becomes empty (cleared), matching spacer row style of row 11.
# Clear values first ...
$ws.Range("B10").Value = $null
$ws.Range("C10").Value = $null
$ws.Range("D10").Value = $null
$ws.Range("E10").Value = $null
# ... then copy only the formatting of E11 (spacer row) onto E10 so its
# style index matches the target (s=12 instead of s=17).
$ws.Range("E11").Copy()
$ws.Range("E10").PasteSpecial(-4122)

# Closing balance line
$ws.Range("D12").Value = "KONTOSTAND AM 10.05.2025"
$ws.Range("E12").Value = "353,17-"

# Next billing date
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 16.05.2025"
